# Uniformity test / tester fix:
# Replace the posthoc-pairwise-test rows (peakValue, RMS, tau, AUC) with the
# Dunn/mannwhitneyu uniformity-test results, reorder each date block so that
# reactionTime, peakTime, difference come first, and drop the now-unused
# "dunn_d" column (K), shrinking the sheet dimension to A1:J22.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-22 (columns A-J)
$data = @(
    @('April16', 'reactionTime', [double]"574.5", [double]"4.731507720099267e-08", $true, [double]"0.6059670781893004", [double]"4.649771442666329e-08", 'Dunn', [double]"-1.046296296296299", 'mannwhitneyu'),
    @('April16', 'peakTime', [double]"1702", [double]"0.122235035799038", $false, [double]"-0.1673525377229081", [double]"0.121469863914346", 'Dunn', [double]"0.7339506172839521", 'mannwhitneyu'),
    @('April16', 'difference', [double]"2104", [double]"6.944883517617799e-05", $true, [double]"-0.4430727023319616", [double]"6.855448672896294e-05", 'Dunn', [double]"1.780246913580246", 'mannwhitneyu'),
    @('April16', 'peakValue', [double]"2881", [double]"2.319259952975573e-18", $true, [double]"-0.9759945130315502", [double]"2.257024464666691e-18", 'Dunn', [double]"61.07962962962966", 'mannwhitneyu'),
    @('April16', 'RMS', [double]"2869", [double]"4.44515293518166e-18", $true, [double]"-0.9677640603566529", [double]"4.32684049725133e-18", 'Dunn', [double]"16.35589506172839", 'mannwhitneyu'),
    @('April16', 'tau', [double]"554", [double]"2.831991584422178e-08", $true, [double]"0.6200274348422496", [double]"2.782633665191067e-08", 'Dunn', [double]"-21.61515123456792", 'mannwhitneyu'),
    @('April16', 'AUC', [double]"2802", [double]"1.518005587749135e-16", $true, [double]"-0.9218106995884774", [double]"1.479448290688163e-16", 'Dunn', [double]"5809.321388888889", 'mannwhitneyu'),
    @('June26', 'reactionTime', [double]"1575.5", [double]"0.4525480764086293", $false, [double]"-0.08058984910836764", [double]"0.4506187202122447", 'Dunn', [double]"0.2765432098765412", 'mannwhitneyu'),
    @('June26', 'peakTime', [double]"1183.5", [double]"0.08026984395540572", $false, [double]"0.1882716049382716", [double]"0.07971977066430727", 'Dunn', [double]"-0.8339506172839535", 'mannwhitneyu'),
    @('June26', 'difference', [double]"1277.5", [double]"0.2659278918781623", $false, [double]"0.1237997256515775", [double]"0.2646022141905987", 'Dunn', [double]"-1.110493827160496", 'mannwhitneyu'),
    @('June26', 'peakValue', [double]"2317", [double]"1.327251541174455e-07", $true, [double]"-0.5891632373113855", [double]"1.30519562593372e-07", 'Dunn', [double]"82.85648148148147", 'mannwhitneyu'),
    @('June26', 'RMS', [double]"2175", [double]"1.070083676788224e-05", $true, [double]"-0.4917695473251029", [double]"1.055027958128037e-05", 'Dunn', [double]"17.32320679012345", 'mannwhitneyu'),
    @('June26', 'tau', [double]"619", [double]"2.57546540999095e-07", $true, [double]"0.575445816186557", [double]"2.533594612822144e-07", 'Dunn', [double]"-42.61480246913578", 'mannwhitneyu'),
    @('June26', 'AUC', [double]"1869", [double]"0.011659015939608", $true, [double]"-0.2818930041152263", [double]"0.01155757319139736", 'Dunn', [double]"4660.227651234569", 'mannwhitneyu'),
    @('May20', 'reactionTime', [double]"1113", [double]"0.03027169445723122", $true, [double]"0.2366255144032922", [double]"0.03003249573877194", 'Dunn', [double]"-0.4185185185185176", 'mannwhitneyu'),
    @('May20', 'peakTime', [double]"916.5", [double]"0.0006033293714184281", $true, [double]"0.3713991769547325", [double]"0.0005963184172090454", 'Dunn', [double]"-1.617901234567899", 'mannwhitneyu'),
    @('May20', 'difference', [double]"1190", [double]"0.09812993066764997", $false, [double]"0.1838134430727023", [double]"0.09750336954318142", 'Dunn', [double]"-1.199382716049378", 'mannwhitneyu'),
    @('May20', 'peakValue', [double]"2794.5", [double]"2.228994531060003e-16", $true, [double]"-0.9166666666666667", [double]"2.172680496951096e-16", 'Dunn', [double]"63.34475308641974", 'mannwhitneyu'),
    @('May20', 'RMS', [double]"2526", [double]"5.408245471040276e-11", $true, [double]"-0.7325102880658436", [double]"5.297948848634982e-11", 'Dunn', [double]"12.87881172839507", 'mannwhitneyu'),
    @('May20', 'tau', [double]"608", [double]"1.791962740654362e-07", $true, [double]"0.5829903978052127", [double]"1.762475032654718e-07", 'Dunn', [double]"-48.5032098765432", 'mannwhitneyu'),
    @('May20', 'AUC', [double]"2040", [double]"0.0003529112438522091", $true, [double]"-0.3991769547325104", [double]"0.0003487918400231757", 'Dunn', [double]"3766.802524691362", 'mannwhitneyu')
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value2 = $row[0]
    $ws.Cells.Item($r, 2).Value2 = $row[1]
    $ws.Cells.Item($r, 3).Value2 = $row[2]
    $ws.Cells.Item($r, 4).Value2 = $row[3]
    $ws.Cells.Item($r, 5).Value2 = $row[4]
    $ws.Cells.Item($r, 6).Value2 = $row[5]
    $ws.Cells.Item($r, 7).Value2 = $row[6]
    $ws.Cells.Item($r, 8).Value2 = $row[7]
    $ws.Cells.Item($r, 9).Value2 = $row[8]
    $ws.Cells.Item($r, 10).Value2 = $row[9]
}

# Remove the obsolete "dunn_d" column (K) entirely; this also shrinks the
# worksheet dimension from A1:K22 down to A1:J22.
$ws.Columns.Item(11).Delete() | Out-Null
